$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "Added plos one to revs, youtube link for brazilian talk, and reorganised
# folders" -- the folder reorganisation swaps the order of the "Biologia" /
# "Pedagogia Musical" block (previously at the bottom, rows 20-21) with the
# "Psicologia" block (previously at rows 6-19), so Biologia/Pedagogia Musical now
# sit right after the Posgrado rows (new rows 6-8) and Psicologia follows (new row 9).
# Four new supervised students are also inserted right under the Psicologia header
# (new rows 10-13), pushing the pre-existing Psicologia student list down to rows 14-26.
# Below, every touched cell is written explicitly so the sheet ends up in that exact
# final layout; rows 1-5 are untouched by the reorganisation and are left alone.

$ws.Range("A6").Value = "Biología"
$ws.Range("B6").Value = "2017 - 2018"
$ws.Range("E6").Value = "Maria Alejandra Abello Mozo  (2018 - 2018)"

$ws.Range("A7").Value = "Pedagogía Musical"
$ws.Range("B7").Value = "2017 - 2019"
$ws.Range("D7").Value = "Universidad Pedagógica Nacional, Colombia"
$ws.Range("E7").Value = "Natalia Elízabeth Moreno Buitrago (2017 - 2019)"

$ws.Range("E8").Value = "Juan Felipe Pérez Ariza (2017 - 2019)"

$ws.Range("A9").Value = "Psicología"
$ws.Range("B9").Value = "Desde 2015"
$ws.Range("C9").Value = "Pregrado"
$ws.Range("D9").Value = "Universidad El Bosque, Colombia"
$ws.Range("E9").Value = "Andrés Castellano Chacón (2017 -2018; supervisión docente, 2019 - actualmente)"

$ws.Range("E10").Value = "Angie Alejandra Lozano Sanjuan (2021 - 2022)"

$ws.Range("E11").Value = "Daniela Martínez Franco (2021 - 2022)"

$ws.Range("E12").Value = "Mariana Saavedra Botero (2021 - 2022)"

$ws.Range("E13").Value = "John Jairo Rubio (2021 - 2022)"

$ws.Range("E14").Value = "Maria Paula Moreno Rodríguez (2019 - 2021)"

$ws.Range("E15").Value = "Andrés Felipe Orozco Serrato (2020 - 2021)"

$ws.Range("E16").Value = "Danny Ferley Gaitan Rodríguez (2019 - 2020)"

$ws.Range("E17").Value = "Hasbleidy Gamboa Ordoñez (2019 - 2020)"

$ws.Range("E18").Value = "Paula Andrea Betancourt Velandia  (2018 - 2019)"

$ws.Range("B19").Value = " "
$ws.Range("E19").Value = "Ana Sofía Gómez Castelblanco (2018 - 2019)"

$ws.Range("E20").Value = "Lina María García Hoyos  (2016 - 2017)"

$ws.Range("E21").Value = "Angie Liliana Pérez Rodríguez  (2016 - 2018)"

$ws.Range("E22").Value = "Lina María Morales Sánchez (2016 - 2017)"

$ws.Range("E23").Value = "Laura Milena Estupiñan Aldana  (2016 - 2017)"

$ws.Range("E24").Value = "Vanesa Díaz Güiza  (2016 - 2018)"

$ws.Range("E25").Value = "Cindy Paola Moncada Gómez (2016 - 2017)"

$ws.Range("E26").Value = "Haydn Ricardo Roldán Morales (2015 - 2016)"

# Clear the cells that held the old (pre-reorganisation) content at these positions
$ws.Range("C6").Value = ""
$ws.Range("B12").Value = ""
$ws.Range("A20").Value = ""
$ws.Range("B20").Value = ""
$ws.Range("D20").Value = ""
$ws.Range("A21").Value = ""
$ws.Range("B21").Value = ""
$ws.Range("D21").Value = ""

# The saved file shows the selection resting on A19 (not the old E22)
[void]$ws.Range("A19").Select()
